# Weekly update of Fruta/Hortaliza prices (Tuna - Agro Chillan)
# Re-sequences Fecha/Volumen/Precio/Unidad/Origen/Precio-Kg columns per row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45041
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("S2").Value = 833
# Row 3
$ws.Range("D3").Value = 45050
$ws.Range("M3").Value = 40
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 14000
$ws.Range("P3").Value = 14000
$ws.Range("S3").Value = 778
# Row 5
$ws.Range("D5").Value = 45096
$ws.Range("N5").Value = 23000
$ws.Range("O5").Value = 23000
$ws.Range("P5").Value = 23000
$ws.Range("Q5").Value = '$/caja 18 kilos'
$ws.Range("R5").Value = 'Región Metropolitana'
$ws.Range("S5").Value = 1278
$ws.Range("T5").Value = 18
# Row 6
$ws.Range("D6").Value = 45014
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 18000
$ws.Range("S6").Value = 1000
# Row 7
$ws.Range("D7").Value = 45020
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 15000
$ws.Range("Q7").Value = '$/caja 16 kilos'
$ws.Range("R7").Value = 'Provincia de Los Andes'
$ws.Range("S7").Value = 938
$ws.Range("T7").Value = 16
# Row 8
$ws.Range("D8").Value = 45001
# Row 9
$ws.Range("D9").Value = 45089
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 22000
$ws.Range("O9").Value = 23000
$ws.Range("P9").Value = 22500
$ws.Range("S9").Value = 1250
# Row 10
$ws.Range("D10").Value = 45091
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 22000
$ws.Range("O10").Value = 22000
$ws.Range("P10").Value = 22000
$ws.Range("R10").Value = 'Región Metropolitana'
$ws.Range("S10").Value = 1222
# Row 11
$ws.Range("D11").Value = 45028
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 18000
$ws.Range("O11").Value = 18000
$ws.Range("P11").Value = 18000
$ws.Range("S11").Value = 1000
# Row 12
$ws.Range("D12").Value = 45002
$ws.Range("M12").Value = 30
$ws.Range("N12").Value = 18000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 18000
$ws.Range("S12").Value = 1000
# Row 13
$ws.Range("D13").Value = 45033
$ws.Range("M13").Value = 60
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 16000
$ws.Range("P13").Value = 15500
$ws.Range("S13").Value = 861
# Row 14
$ws.Range("D14").Value = 45049
$ws.Range("M14").Value = 80
$ws.Range("N14").Value = 15000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 15000
$ws.Range("S14").Value = 833
# Row 15
$ws.Range("D15").Value = 45099
$ws.Range("M15").Value = 40
$ws.Range("N15").Value = 22000
$ws.Range("O15").Value = 22000
$ws.Range("P15").Value = 22000
$ws.Range("S15").Value = 1222
# Row 16
$ws.Range("D16").Value = 45044
$ws.Range("M16").Value = 60
$ws.Range("N16").Value = 15000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 15000
$ws.Range("S16").Value = 833
# Row 17
$ws.Range("D17").Value = 44999
$ws.Range("N17").Value = 17000
$ws.Range("O17").Value = 18000
$ws.Range("P17").Value = 17500
$ws.Range("S17").Value = 972
# Row 18
$ws.Range("D18").Value = 45062
$ws.Range("M18").Value = 90
$ws.Range("N18").Value = 13000
$ws.Range("O18").Value = 14000
$ws.Range("P18").Value = 13444
$ws.Range("S18").Value = 747
# Row 19
$ws.Range("D19").Value = 45030
$ws.Range("M19").Value = 40
$ws.Range("N19").Value = 18000
$ws.Range("O19").Value = 18000
$ws.Range("P19").Value = 18000
$ws.Range("S19").Value = 1000
# Row 20
$ws.Range("D20").Value = 45021
$ws.Range("M20").Value = 60
$ws.Range("N20").Value = 15000
$ws.Range("O20").Value = 16000
$ws.Range("P20").Value = 15500
$ws.Range("R20").Value = 'Provincia de Los Andes'
$ws.Range("S20").Value = 861
# Row 21
$ws.Range("D21").Value = 45043
$ws.Range("N21").Value = 15000
$ws.Range("O21").Value = 15000
$ws.Range("P21").Value = 15000
$ws.Range("S21").Value = 833
# Row 22
$ws.Range("D22").Value = 45037
$ws.Range("N22").Value = 16000
$ws.Range("O22").Value = 16000
$ws.Range("P22").Value = 16000
$ws.Range("S22").Value = 889
